# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain stored as TEXT
# (matches the original inlineStr/shared-string type in the target workbook).
# Setting NumberFormat to "@" (Text) before assigning the value keeps it a string
# instead of Excel auto-converting it to a numeric value.
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D12", "D16", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D29", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "57.098.97"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "2.418.10"
$ws.Range("E3").Value = "  -3.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "488.66"
$ws.Range("E5").Value = "  -1.14%  "

# Row 6
$ws.Range("D6").Value = "154.56"
$ws.Range("E6").Value = "  +1.05%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  +19.10%  "

# Row 9
$ws.Range("D9").Value = "2.445.59"
$ws.Range("E9").Value = "  -3.08%  "

# Row 10
$ws.Range("D10").Value = "6.35"
$ws.Range("E10").Value = "  +10.12%  "

# Row 11
$ws.Range("D11").Value = "0.0999"
$ws.Range("E11").Value = "  -0.86%  "

# Row 12
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13
$ws.Range("E13").Value = "  +1.44%  "

# Row 14
$ws.Range("D14").Value = "2.842.34"
$ws.Range("E14").Value = "  -3.40%  "

# Row 15
$ws.Range("D15").Value = "57.095.65"
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").Value = "20.63"
$ws.Range("E16").Value = "  -3.21%  "

# Row 17
$ws.Range("E17").Value = "  -3.72%  "

# Row 18
$ws.Range("D18").Value = "2.437.62"
$ws.Range("E18").Value = "  -4.02%  "

# Row 19
$ws.Range("D19").Value = "4.69"
$ws.Range("E19").Value = "  +2.09%  "

# Row 20
$ws.Range("D20").Value = "324.55"
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("D21").Value = "10.01"
$ws.Range("E21").Value = "  -3.07%  "

# Row 22
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").Value = "57.79"
$ws.Range("E24").Value = "  -1.15%  "

# Row 25
$ws.Range("E25").Value = "  -1.54%  "

# Row 26
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.41%  "

# Row 27
$ws.Range("D27").Value = "0.162"
$ws.Range("E27").Value = "  -1.23%  "

# Row 28
$ws.Range("D28").Value = "2.526.74"
$ws.Range("E28").Value = "  -3.55%  "

# Row 29
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  -4.52%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  -5.24%  "

# Row 31
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("D32").Value = "151.06"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").Value = "18.68"
$ws.Range("E33").Value = "  +1.79%  "

# Row 34
$ws.Range("E34").Value = "  -0.49%  "

# Row 35
$ws.Range("D35").Value = "5.29"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  -1.09%  "

# Row 38
$ws.Range("D38").Value = "0.831"
$ws.Range("E38").Value = "  -6.55%  "

# Row 39
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +7.98%  "

# Row 40
$ws.Range("D40").Value = "34.06"

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  -0.46%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.37"
$ws.Range("E42").Value = "  -2.65%  "

# Row 43
$ws.Range("D43").Value = "279.54"
$ws.Range("E43").Value = "  +3.99%  "

# Row 44
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "0.599"
$ws.Range("E45").Value = "  -3.50%  "

# Row 46
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  -5.26%  "

# Row 47
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("D48").Value = "0.0227"
$ws.Range("E48").Value = "  -1.34%  "

# Row 49
$ws.Range("D49").Value = "4.59"
$ws.Range("E49").Value = "  -7.10%  "

# Row 50
$ws.Range("D50").Value = "1.901.80"
$ws.Range("E50").Value = "  +0.10%  "

# Row 51
$ws.Range("D51").Value = "17.59"
$ws.Range("E51").Value = "  -2.32%  "

